# Weekly update: insert a new price record at row 28 for
# "Feria Lagunitas de Puerto Montt" / Poroto verde, shifting the
# existing historical rows down by one, and append the row that falls
# off the bottom of the used range (previous row 120) as the new last
# row (121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28. This shifts rows 28:120
# down to 29:121 (carrying all their values/styles with them) and
# extends the worksheet's used range to A1:R121.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with this week's record. Most of
# the descriptive columns are identical to the template row that used
# to live at row 28 (now at row 29): market, region, category,
# variety, quality, unit of sale, origin, and classification.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44953
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Magnum"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 30000
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 1200
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
